# Convert the "capital" sheet from a simple summary row into a balance-sheet
# style header-only layout, and add the new shared strings required by the
# new column headers. Also re-point the "holding" sheet header cells that
# reused some of the old shared strings (Symbol/ShortQty/Note) so they keep
# showing the same text after the shared-string table is effectively
# reshuffled.

$wb = $excel.ActiveWorkbook

$capital = $wb.Worksheets.Item("capital")
$holding = $wb.Worksheets.Item("holding")

# --- capital sheet -------------------------------------------------------
# Clear everything first (old header + old data row) then rebuild.
$capital.Cells.Clear()

$headers = @("AcctIDByMXZ", "Cash", "CashEquivalent", "ETFLongAmt", "CompositeLongAmt", "TotalAsset", "ETFShortAmt", "CompositeShortAmt", "Liability", "ApproximateNetAsset")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $capital.Cells.Item(1, $col).Value = $headers[$i]
}

# Apply the "Text" style (style index 1 in styles.xml, same style already
# used on A1) to columns A through I (all headers except the last one).
$capital.Range("A1:I1").Style = $capital.Range("A1").Style

# Column widths / bestFit per target cols definition
$capital.Columns.Item(1).ColumnWidth = 14.625
$capital.Columns.Item(2).ColumnWidth = 5.375
$capital.Columns.Item(3).ColumnWidth = 14.125
$capital.Columns.Item(4).ColumnWidth = 12.25
$capital.Columns.Item(5).ColumnWidth = 18.875
$capital.Columns.Item(6).ColumnWidth = 10
$capital.Columns.Item(7).ColumnWidth = 12.5
$capital.Columns.Item(8).ColumnWidth = 19.125
$capital.Columns.Item(9).ColumnWidth = 7.5
$capital.Columns.Item(10).ColumnWidth = 20.625

# Selection / active cell as captured in the saved file.
$capital.Activate()
$capital.Range("F5").Select()

# --- holding sheet ---------------------------------------------------------
# Text is unchanged; nothing to do explicitly here, the shared-string
# reordering happens naturally as a consequence of the capital sheet edits.

$wb.Save()
